# Deploying to gh-pages from @ Alvearie/alvearie-fhir-ig@8e4a450c507ef6f746e072652acbb72e9504f19a
#
# Updates the ValueSet-quality-measure-population-type workbook:
#  - Metadata sheet: bump Version 5.0.0 -> 6.0.0, refresh Date, fill in the
#    Publisher value, drop the duplicated "Contact" row and replace it with
#    a Jurisdiction row.
#  - Exclude-from-MeasurePopulation sheet: fix the "exlusion" typo to
#    "exclusion".

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Metadata")

# The sheet currently has two identical "Contact" / "No display for
# ContactDetail" rows (rows 10 and 11). Remove the duplicate (row 11); this
# shifts every row below it up by one, matching the new A1:B14 dimension.
$ws.Rows("11").Delete()

# Version
$ws.Range("B3").Value = "6.0.0"

# Date
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank before; now populated.
$ws.Range("B9").Value = "Alvearie Team"

# The remaining Contact / No-display-for-ContactDetail row becomes a
# Jurisdiction row.
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Fix the "measure-population-exlusion" typo on the exclusion sheet.
$ws4 = $wb.Worksheets.Item("Exclude from MeasurePopulatio")
$ws4.Range("A3").Value = "measure-population-exclusion"
